$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N2").Value = "2018-12-31 00:00:00"
$ws.Range("O2").Value = 1704774114.27
$ws.Range("P2").Value = 345508618.63
$ws.Range("Q2").Value = 120234313.34
$ws.Range("R2").Value = 116.4749780426
$ws.Range("S2").Value = 212299924.22
$ws.Range("T2").Value = 0.0854792858
$ws.Range("U2").Value = 375529711.28
$ws.Range("V2").Value = 3.517392671
$ws.Range("W2").Value = 799311478.38
$ws.Range("X2").Value = 230697347.59
$ws.Range("Y2").Value = -6.4140407448
$ws.Range("Z2").Value = 2828033.16
$ws.Range("AA2").Value = -86.7788627394
$ws.Range("AB2").Value = 905462635.89
$ws.Range("AC2").Value = 15.0604885437
$ws.Range("AD2").Value = 13.4727646257
$ws.Range("AE2").Value = 11.7263044316
$ws.Range("AF2").Value = 157.1592200639
$ws.Range("AG2").Value = 46.8866503597
